# Resume/CV update ("updated GPA in resume and cv").
#
# Substantive change: GPA 3.5 -> 3.6 in both the M.Sc. and the B.Sc.
# education rows.
#
# The same save also shows the "plotly" and "Keras" skill entries losing
# their own run / stale spell-check markup and folding back into the
# text that precedes them (", plotly" and " Pytorch, Keras" each become
# a single run) - same wording, just tidied up. We reproduce that with
# an in-place Find/Replace across those spans, which merges the runs and
# drops the now-stale proofing marks.
#
# Find.Execute's Wrap:=wdFindContinue (1) makes a single call walk the
# whole story and replace every match (both GPA rows included), so one
# call per phrase is enough.

$d = $word.ActiveDocument

# 1) Data Visualization row: ", " + "plotly" -> single ", plotly" run.
$found1 = $d.Content.Find.Execute(", plotly", $true, $false, $false, $false, `
    $false, $true, 1, $false, ", plotly", 2)
Write-Output "ggplot2/plotly run tidy-up found: $found1"

# 2) Deep Learning row: " Pytorch, " + "Keras" -> single " Pytorch, Keras" run.
$found2 = $d.Content.Find.Execute(" Pytorch, Keras", $true, $false, $false, $false, `
    $false, $true, 1, $false, " Pytorch, Keras", 2)
Write-Output "Pytorch/Keras run tidy-up found: $found2"

# 3) GPA: 3.5 -> GPA: 3.6 (M.Sc. row and B.Sc. row - both updated by this
#    single wrapping Find/Replace call).
$found3 = $d.Content.Find.Execute("GPA: 3.5", $true, $false, $false, $false, `
    $false, $true, 1, $false, "GPA: 3.6", 2)
Write-Output "GPA 3.5 -> 3.6 found: $found3"
